$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the student-ID values up by one row for rows 5-7 (A5,A6,A7) and
# correspondingly update only the leading "Student ID:  <id>" line of the
# matching B cell, leaving the rest of each B cell's text (URL/questions)
# untouched. Then delete the now-duplicate last row (row 8).

function Set-NewId($row, $newId, $oldId) {
    $aCell = $ws.Cells.Item($row, 1)
    $bCell = $ws.Cells.Item($row, 2)

    $aCell.Value = $newId

    $bText = $bCell.Value2
    $lines = $bText -split "`n", 2
    $lines[0] = $lines[0].Replace($oldId, $newId)
    $bCell.Value = ($lines -join "`n")
}

Set-NewId 5 "q0762379" "q0328135_previewuser"
Set-NewId 6 "q1371623" "q0762379"
Set-NewId 7 "q1411379" "q1371623"

# Remove the old row 8 (now a duplicate of row 7's content).
$ws.Rows.Item(8).Delete()
